{"js": "// Apply the \"Compact\" paragraph style to the title paragraph\n// (\"An\u00e1lisis Multivariante\"), which is the first paragraph of the\n// document body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Prefer matching on the known text so the script is resilient even if\n// paragraph order ever shifts; fall back to the very first paragraph\n// (which is where the title lives in this document).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text && p.text.trim() === \"An\u00e1lisis Multivariante\") {\n    target = p;\n    break;\n  }\n}\nif (!target) {\n  target = paragraphs.items[0];\n}\n\ntarget.style = \"Compact\";\nawait context.sync();\n", "ps1": "# Apply the \"Compact\" paragraph style to the title paragraph\n# (\"An\u00e1lisis Multivariante\"), which is the first paragraph of the\n# document body.\n$d = $word.ActiveDocument\n\n$target = $null\n\n# Prefer matching on the known title text so the script is resilient even\n# if paragraph order ever shifts.\n$range = $d.Content\nif ($range.Find.Execute(\"An\u00e1lisis Multivariante\")) {\n    $target = $range.Paragraphs(1)\n}\n\n# Fall back to the very first paragraph of the document (where the title\n# lives in this document) if the text search didn't resolve.\nif (-not $target) {\n    $target = $d.Paragraphs(1)\n}\n\n$target.Style = \"Compact\"\n"}
